# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values on the zh-cn and de-de
# report sheets to reflect a newer report-generation run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-14 08:15:25"
$zhcn.Range("H4").Value = "2016-03-14 08:15:43"
$zhcn.Range("E5").Value = "2016-03-14 08:15:25"
$zhcn.Range("H5").Value = "2016-03-14 08:15:43"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-14 08:15:29"
$dede.Range("H4").Value = "2016-03-14 08:15:49"
$dede.Range("E5").Value = "2016-03-14 08:15:29"
$dede.Range("H5").Value = "2016-03-14 08:15:49"
